$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" '96.092.66'
Set-TextCell $ws "E2" '  +4.51%  '
Set-TextCell $ws "D3" '3.672.16'
Set-TextCell $ws "E3" '  +10.28%  '
Set-TextCell $ws "E4" '  +0.16%  '
Set-TextCell $ws "D5" '241.52'
Set-TextCell $ws "D6" '643.44'
Set-TextCell $ws "E6" '  +4.64%  '
Set-TextCell $ws "E7" '  +5.06%  '
Set-TextCell $ws "D8" '0.402'
Set-TextCell $ws "E8" '  +4.79%  '
Set-TextCell $ws "D9" '1.00'
Set-TextCell $ws "E9" '  -0.02%  '
Set-TextCell $ws "E10" '  +5.23%  '
Set-TextCell $ws "D11" '3.676.36'
Set-TextCell $ws "E11" '  +10.40%  '
Set-TextCell $ws "D12" '43.93'
Set-TextCell $ws "E12" '  +2.45%  '
Set-TextCell $ws "E13" '  +3.79%  '
Set-TextCell $ws "D14" '6.39'
Set-TextCell $ws "E14" '  +3.72%  '
Set-TextCell $ws "D15" '4.366.58'
Set-TextCell $ws "E15" '  +10.50%  '
Set-TextCell $ws "D16" '95.986.64'
Set-TextCell $ws "E16" '  +4.61%  '
Set-TextCell $ws "E17" '  +5.69%  '
Set-TextCell $ws "B18" 'Uniswap'
Set-TextCell $ws "C18" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws "D18" '13.64'
Set-TextCell $ws "E18" '  +25.76%  '
Set-TextCell $ws "B19" 'WrappedEther'
Set-TextCell $ws "C19" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws "D19" '3.661.04'
Set-TextCell $ws "E19" '  +9.78%  '
Set-TextCell $ws "B20" 'Polkadot'
Set-TextCell $ws "C20" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws "D20" '8.03'
Set-TextCell $ws "E20" '  -0.59%  '
Set-TextCell $ws "D21" '18.86'
Set-TextCell $ws "E21" '  +8.88%  '
Set-TextCell $ws "D22" '521.03'
Set-TextCell $ws "E22" '  +5.59%  '
Set-TextCell $ws "D23" '3.44'
Set-TextCell $ws "E23" '  +0.66%  '
Set-TextCell $ws "D24" '0.483'
Set-TextCell $ws "E24" '  +8.96%  '
Set-TextCell $ws "E25" '  +8.91%  '
Set-TextCell $ws "D26" '6.82'
Set-TextCell $ws "E26" '  +5.85%  '
Set-TextCell $ws "D27" '97.80'
Set-TextCell $ws "E27" '  +8.62%  '
Set-TextCell $ws "D28" '12.62'
Set-TextCell $ws "E28" '  +5.84%  '
Set-TextCell $ws "D29" '3.20'
Set-TextCell $ws "E29" '  +23.05%  '
Set-TextCell $ws "D30" '11.74'
Set-TextCell $ws "E30" '  +5.26%  '
Set-TextCell $ws "E31" '  +2.32%  '
Set-TextCell $ws "D32" '1.00'
Set-TextCell $ws "E32" '  -0.07%  '
Set-TextCell $ws "B33" 'EthereumClassic'
Set-TextCell $ws "C33" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws "D33" '33.07'
Set-TextCell $ws "E33" '  +16.64%  '
Set-TextCell $ws "B34" 'Cronos'
Set-TextCell $ws "C34" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws "D34" '0.181'
Set-TextCell $ws "E34" '  +4.99%  '
Set-TextCell $ws "B35" 'Binance-PegBSC-USD'
Set-TextCell $ws "C35" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws "D35" '0.998'
Set-TextCell $ws "E35" '  -1.29%  '
Set-TextCell $ws "D36" '0.581'
Set-TextCell $ws "E36" '  +10.04%  '
Set-TextCell $ws "D37" '566.52'
Set-TextCell $ws "E37" '  -0.18%  '
Set-TextCell $ws "B38" 'Fetch.AI'
Set-TextCell $ws "C38" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws "D38" '1.50'
Set-TextCell $ws "E38" '  +9.97%  '
Set-TextCell $ws "B39" 'RenderToken'
Set-TextCell $ws "C39" 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextCell $ws "D39" '7.88'
Set-TextCell $ws "E39" '  +6.83%  '
Set-TextCell $ws "D40" '0.966'
Set-TextCell $ws "E40" '  +11.32%  '
Set-TextCell $ws "E41" '  +3.23%  '
Set-TextCell $ws "E42" '  -0.07%  '
Set-TextCell $ws "B43" 'Filecoin'
Set-TextCell $ws "C43" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws "D43" '5.82'
Set-TextCell $ws "E43" '  +7.51%  '
Set-TextCell $ws "B44" 'ImmutableX'
Set-TextCell $ws "C44" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws "D44" '1.75'
Set-TextCell $ws "E44" '  +4.49%  '
Set-TextCell $ws "B45" 'VeChain'
Set-TextCell $ws "C45" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws "D45" '0.0432'
Set-TextCell $ws "E45" '  +4.62%  '
Set-TextCell $ws "D46" '23.75'
Set-TextCell $ws "E46" '  +0.25%  '
Set-TextCell $ws "B47" 'EnergySwap'
Set-TextCell $ws "C47" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws "D47" '33.74'
Set-TextCell $ws "E47" '  +50.11%  '
Set-TextCell $ws "B48" 'Stacks'
Set-TextCell $ws "C48" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws "D48" '2.23'
Set-TextCell $ws "E48" '  +5.91%  '
Set-TextCell $ws "B49" 'Cosmos'
Set-TextCell $ws "C49" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws "D49" '8.34'
Set-TextCell $ws "E49" '  +4.44%  '
Set-TextCell $ws "B50" 'OKB'
Set-TextCell $ws "C50" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws "D50" '54.27'
Set-TextCell $ws "E50" '  +4.74%  '
Set-TextCell $ws "B51" 'MantraDAO'
Set-TextCell $ws "C51" 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextCell $ws "D51" '3.50'
Set-TextCell $ws "E51" '  -2.66%  '
